$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking Price values to stay text (source data is inline-string formatted),
# matching the original inlineStr content (e.g. "317.64", "1.240") rather than being
# auto-converted to numbers by Excel.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.683.05'
$ws.Range("E2").Value = '  +0.93%  '

$ws.Range("D3").Value = '1.804.78'
$ws.Range("E3").Value = '  -0.53%  '

$ws.Range("E4").Value = '  +0.52%  '

$ws.Range("D5").Value = '317.64'
$ws.Range("E5").Value = '  -0.12%  '

$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("D7").Value = '0.5343'
$ws.Range("E7").Value = '  -5.70%  '

$ws.Range("D8").Value = '0.3784'
$ws.Range("E8").Value = '  -1.41%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.07523'
$ws.Range("E9").Value = '  -1.44%  '

$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").Value = '42.47'
$ws.Range("E10").Value = '  -2.29%  '

$ws.Range("E11").Value = '  -1.87%  '

$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.60%  '

$ws.Range("D13").Value = '20.83'
$ws.Range("E13").Value = '  -2.45%  '

$ws.Range("D14").Value = '6.182'
$ws.Range("E14").Value = '  -0.92%  '

$ws.Range("D15").Value = '7.369'
$ws.Range("E15").Value = '  +1.67%  '

$ws.Range("D16").Value = '1.803.76'
$ws.Range("E16").Value = '  -0.12%  '

$ws.Range("E17").Value = '  -1.89%  '

$ws.Range("D18").Value = '0.00001068'
$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("D19").Value = '0.06454'
$ws.Range("E19").Value = '  -0.96%  '

$ws.Range("E20").Value = '  +0.41%  '

$ws.Range("D21").Value = '17.27'
$ws.Range("E21").Value = '  -0.09%  '

$ws.Range("D22").Value = '5.917'
$ws.Range("E22").Value = '  -1.41%  '

$ws.Range("D23").Value = '28.704.64'
$ws.Range("E23").Value = '  +0.99%  '

$ws.Range("E24").Value = '  -1.52%  '

$ws.Range("D25").Value = '2.105'
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("D26").Value = '160.35'
$ws.Range("E26").Value = '  +2.25%  '

$ws.Range("D27").Value = '20.46'
$ws.Range("E27").Value = '  -1.74%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '2.383'
$ws.Range("E28").Value = '  -0.29%  '

$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").Value = '2.014.85'
$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").Value = '123.22'
$ws.Range("E30").Value = '  -0.37%  '

$ws.Range("D31").Value = '1.104'
$ws.Range("E31").Value = '  -3.77%  '

$ws.Range("D32").Value = '0.1053'
$ws.Range("E32").Value = '  +0.38%  '

$ws.Range("D33").Value = '5.657'
$ws.Range("E33").Value = '  -2.02%  '

$ws.Range("D34").Value = '3.675'
$ws.Range("E34").Value = '  +1.36%  '

$ws.Range("D35").Value = '0.2261'
$ws.Range("E35").Value = '  +5.73%  '

$ws.Range("E36").Value = '  +5.97%  '

$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = '8.893'
$ws.Range("E37").Value = '  +1.96%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02308'
$ws.Range("E38").Value = '  -0.45%  '

$ws.Range("D39").Value = '1.240'
$ws.Range("E39").Value = '  +7.39%  '

$ws.Range("D40").Value = '5.047'
$ws.Range("E40").Value = '  +0.05%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '11.30'
$ws.Range("E41").Value = '  -3.21%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.6252'
$ws.Range("E42").Value = '  -2.82%  '

$ws.Range("E43").Value = '  +0.43%  '

$ws.Range("E44").Value = '  +1.21%  '

$ws.Range("D45").Value = '13.30'
$ws.Range("E45").Value = '  -1.49%  '

$ws.Range("D46").Value = '0.5880'
$ws.Range("E46").Value = '  -2.08%  '

$ws.Range("D47").Value = '3.699'

$ws.Range("D48").Value = '126.14'
$ws.Range("E48").Value = '  +2.95%  '

$ws.Range("D49").Value = '1.953'
$ws.Range("E49").Value = '  +0.70%  '

$ws.Range("E50").Value = '  +0.57%  '

$ws.Range("D51").Value = '0.06885'
$ws.Range("E51").Value = '  +0.67%  '
